$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RR")

# Column E ("FHIR" mapping) gains new text for rows 66-78, mirroring the
# style already used by the existing E64/E65 cells (Responsible Agency
# address/telecom). Row heights grow to fit the wrapped text.
#
# This rounds out the Responsible Agency extension (text/identifier/name)
# and adds full FHIR-path coverage for the Routing Entity and
# Rules Authoring Agency extensions.

$rowData = @(
    @{ Row = 66; Height = 72; Text = "rr-communication (Communication).payload:sliceRelevantReportableConditionInformation.rr-plandefinition (PlanDefinition).jurisdiction.extension-rr-responsible-agency.rr-responsible-agency (US Core Organization).text" },
    @{ Row = 67; Height = 72; Text = "rr-communication (Communication).payload:sliceRelevantReportableConditionInformation.rr-plandefinition (PlanDefinition).jurisdiction.extension-rr-responsible-agency.rr-responsible-agency (US Core Organization).identifier" },
    @{ Row = 68; Height = 72; Text = "rr-communication (Communication).payload:sliceRelevantReportableConditionInformation.rr-plandefinition (PlanDefinition).jurisdiction.extension-rr-responsible-agency.rr-responsible-agency (US Core Organization).name" },
    @{ Row = 69; Height = 72; Text = "rr-communication (Communication).payload:sliceRelevantReportableConditionInformation.rr-plandefinition (PlanDefinition).jurisdiction.extension-rr-routing-entity.rr-routing-entity (US Core Organization).address" },
    @{ Row = 70; Height = 72; Text = "rr-communication (Communication).payload:sliceRelevantReportableConditionInformation.rr-plandefinition (PlanDefinition).jurisdiction.extension-rr-routing-entity.rr-routing-entity (US Core Organization).telecom" },
    @{ Row = 71; Height = 60; Text = "rr-communication (Communication).payload:sliceRelevantReportableConditionInformation.rr-plandefinition (PlanDefinition).jurisdiction.extension-rr-routing-entity.rr-routing-entity (US Core Organization).text" },
    @{ Row = 72; Height = 72; Text = "rr-communication (Communication).payload:sliceRelevantReportableConditionInformation.rr-plandefinition (PlanDefinition).jurisdiction.extension-rr-routing-entity.rr-routing-entity (US Core Organization).identifier" },
    @{ Row = 73; Height = 60; Text = "rr-communication (Communication).payload:sliceRelevantReportableConditionInformation.rr-plandefinition (PlanDefinition).jurisdiction.extension-rr-routing-entity.rr-routing-entity (US Core Organization).name" },
    @{ Row = 74; Height = 72; Text = "rr-communication (Communication).payload:sliceRelevantReportableConditionInformation.rr-plandefinition (PlanDefinition).jurisdiction.extension-rr-rules-authoring-agency.rr-rules-authoring-agency (US Core Organization).address" },
    @{ Row = 75; Height = 72; Text = "rr-communication (Communication).payload:sliceRelevantReportableConditionInformation.rr-plandefinition (PlanDefinition).jurisdiction.extension-rr-rules-authoring-agency.rr-rules-authoring-agency (US Core Organization).telecom" },
    @{ Row = 76; Height = 60; Text = "rr-communication (Communication).payload:sliceRelevantReportableConditionInformation.rr-plandefinition (PlanDefinition).jurisdiction.extension-rr-rules-authoring-agency.rr-rules-authoring-agency (US Core Organization).text" },
    @{ Row = 77; Height = 72; Text = "rr-communication (Communication).payload:sliceRelevantReportableConditionInformation.rr-plandefinition (PlanDefinition).jurisdiction.extension-rr-rules-authoring-agency.rr-rules-authoring-agency (US Core Organization).identifier" },
    @{ Row = 78; Height = 60; Text = "rr-communication (Communication).payload:sliceRelevantReportableConditionInformation.rr-plandefinition (PlanDefinition).jurisdiction.extension-rr-rules-authoring-agency.rr-rules-authoring-agency (US Core Organization).name" }
)

$formatSource = $ws.Cells.Item(64, 5)
$formatSource.Copy()

foreach ($item in $rowData) {
    $row = $item.Row
    $ws.Rows.Item($row).RowHeight = $item.Height
    $cell = $ws.Cells.Item($row, 5)
    $cell.PasteSpecial(-4122)
    $cell.Value2 = $item.Text
}

$excel.CutCopyMode = $false

# Keep the sheet's remembered selection in sync with the new bottom row.
$ws.Range("C77").Select()
